$d = $word.ActiveDocument

$replacements = @(
    @{old = "55×16="; new = "42×37="},
    @{old = "16×14="; new = "86×34="},
    @{old = "19×12="; new = "12×39="},
    @{old = "55×77="; new = "64×28="},
    @{old = "31×42="; new = "41×80="},
    @{old = "73×25="; new = "14×16="},
    @{old = "42×49="; new = "79×57="},
    @{old = "37×67="; new = "25×55="},
    @{old = "86×60="; new = "72×60="},
    @{old = "25×66="; new = "42×31="},
    @{old = "80×19="; new = "34×21="},
    @{old = "19×77="; new = "67×96="},
    @{old = "21×65="; new = "79×88="},
    @{old = "30×60="; new = "99×59="},
    @{old = "71×95="; new = "85×13="},
    @{old = "14×65="; new = "50×57="},
    @{old = "66×63="; new = "64×22="},
    @{old = "85×79="; new = "81×91="},
    @{old = "99×55="; new = "71×67="},
    @{old = "49×63="; new = "53×38="},
    @{old = "77×72="; new = "33×53="},
    @{old = "47×70="; new = "69×67="},
    @{old = "89×18="; new = "95×47="},
    @{old = "27×16="; new = "11×59="},
    @{old = "43×30="; new = "64×51="}
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}

$d.Save()
